# Update cryptos list data (Price and Volume(1h) columns) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.367.25'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '3.101.95'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '3.093.07'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.12%  '
$ws.Range('E11').Value = '  +5.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '3.606.43'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '63.455.96'
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.111'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '3.109.61'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '503.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.70'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.708'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.38'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.04'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.53'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.75%  '
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.38'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '534.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.91'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.18'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0415'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0798'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('D40').Value = '3.083.28'
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.118'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.12'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.256'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.09'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.33%  '
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '0.0₃0500'
$ws.Range('E50').Value = '  -5.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +63.99%  '
